$d = $word.ActiveDocument

# Namespace-qualified "mini package" wrapper so InsertXML gets a full OOXML
# part (lets us control the exact run/tab/bookmark structure being inserted,
# instead of letting plain text edits re-flow/merge the surrounding runs).
function New-PkgXml([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- 1. Drop the referring physician's name ("Dra. Katalina Parra") -------
# After the edit only "Referencia:" + a tab should remain on that line.
$found = $d.Content
$null = $found.Find.Execute(":	Dra. Katalina Parra")
$refRange = $d.Range($found.Start, $found.End)

$refXml = New-PkgXml(
    '<w:p>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr><w:t>:</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr><w:tab/></w:r>' +
    '</w:p>'
)
$refRange.InsertXML($refXml)

# --- 2. Relocate the "_GoBack" bookmark ------------------------------------
# Remove it from its old spot (a blank paragraph near the end of the doc) ...
try {
    $d.Bookmarks.Item("_GoBack").Delete()
} catch {
    # not present - nothing to clean up
}

# ...and re-insert it mid-sentence in the HISTORIA paragraph, right before
# "actualmente se encuentra en el tercer mes de terapia de lenguaje."
$found2 = $d.Content
$null = $found2.Find.Execute("actualmente se encuentra en el tercer mes de terapia de lenguaje.")
$histRange = $d.Range($found2.Start, $found2.End)

$histXml = New-PkgXml(
    '<w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Book Antiqua" w:hAnsi="Book Antiqua"/><w:sz w:val="20"/><w:szCs w:val="22"/></w:rPr>' +
    '<w:t>actualmente se encuentra en el tercer mes de terapia de lenguaje.</w:t></w:r></w:p>'
)
$histRange.InsertXML($histXml)
